$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Range("Q2").Value = "51539157"
$ws.Range("Q3").Value = "51539148"
$ws.Range("R3").Value = "51539149"
$ws.Range("AD3").Value = "05-11-2022"
$ws.Range("Q4").Value = "51539150"
